# Admin Panel "Change Password" feature -- data sheet updates.
# (The app/UI code changes described in the commit message live outside
# this workbook; here we apply the resulting worksheet data changes.)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# weekend_tasks: add a "Goal Credits" column per kid, add new goal rows.
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("weekend_tasks")

# Header row
$ws.Cells.Item(1,1).Value = "Jackson"
$ws.Cells.Item(1,2).Value = "Jackson Goal"
$ws.Cells.Item(1,3).Value = "Jackson Goal Credits"
$ws.Cells.Item(1,4).Value = "Natalie"
$ws.Cells.Item(1,5).Value = "Natalie Goal"
$ws.Cells.Item(1,6).Value = "Natalie Goal Credits"
$ws.Cells.Item(1,7).Value = "Brooke"
$ws.Cells.Item(1,8).Value = "Brooke Goal"
$ws.Cells.Item(1,9).Value = "Brooke Goal Credits"

# Row 2
$ws.Cells.Item(2,1).Value = "play baseball game"
$ws.Cells.Item(2,2).Value = "Hit Baseball Over Infeed"
$ws.Cells.Item(2,3).Value = 5
$ws.Cells.Item(2,4).Value = "play softball game"
$ws.Cells.Item(2,5).Value = "Get a Out"
$ws.Cells.Item(2,6).Value = 5
$ws.Cells.Item(2,7).Value = "play softball game"
$ws.Cells.Item(2,8).Value = "Hit a Pitch from a Kid"
$ws.Cells.Item(2,9).Value = 5

# Row 3 (new)
$ws.Cells.Item(3,1).Value = ""
$ws.Cells.Item(3,2).Value = "Get a Out"
$ws.Cells.Item(3,3).Value = 5
$ws.Cells.Item(3,4).Value = ""
$ws.Cells.Item(3,5).Value = "Hit a Double"
$ws.Cells.Item(3,6).Value = 5
$ws.Cells.Item(3,7).Value = ""
$ws.Cells.Item(3,8).Value = "Get a Out"
$ws.Cells.Item(3,9).Value = 5

# Row 4 (new)
$ws.Cells.Item(4,1).Value = ""
$ws.Cells.Item(4,2).Value = ""
$ws.Cells.Item(4,3).Value = ""
$ws.Cells.Item(4,4).Value = ""
$ws.Cells.Item(4,5).Value = "Go a Whole Week Without Yelling at Mom"
$ws.Cells.Item(4,6).Value = 25
$ws.Cells.Item(4,7).Value = ""
$ws.Cells.Item(4,8).Value = ""
$ws.Cells.Item(4,9).Value = ""

# ---------------------------------------------------------------------
# morning_options: refresh Special Breakfast / Snacks options, and add
# two new people/rows at the bottom (Keish, Breakfast Oats).
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("morning_options")

$ws2.Cells.Item(2,3).Value = "Doughnuts"
$ws2.Cells.Item(2,11).Value = "Fruit Rolls"

$ws2.Cells.Item(3,3).Value = "French Toast"
$ws2.Cells.Item(3,11).Value = "Granola bar"

$ws2.Cells.Item(4,3).Value = "Pancakes"
$ws2.Cells.Item(4,11).Value = "Apple Sauce Pouch"

$ws2.Cells.Item(5,3).Value = "Waffles"
$ws2.Cells.Item(5,11).Value = "Pretzels"

$ws2.Cells.Item(6,3).Value = ""
$ws2.Cells.Item(6,4).Value = ""
$ws2.Cells.Item(6,11).Value = "Ritz Crackers"

$ws2.Cells.Item(7,11).Value = "Doritos"
$ws2.Cells.Item(7,12).Value = "yes"

$ws2.Cells.Item(8,11).Value = "Sun Chips"
$ws2.Cells.Item(8,12).Value = "yes"

# New row 15
$ws2.Cells.Item(15,1).Value = "Keish"
$ws2.Cells.Item(15,2).Value = "yes"
$ws2.Cells.Item(15,3).Value = ""
$ws2.Cells.Item(15,4).Value = ""
$ws2.Cells.Item(15,5).Value = ""
$ws2.Cells.Item(15,6).Value = ""
$ws2.Cells.Item(15,7).Value = ""
$ws2.Cells.Item(15,8).Value = ""
$ws2.Cells.Item(15,9).Value = ""
$ws2.Cells.Item(15,10).Value = ""
$ws2.Cells.Item(15,11).Value = ""
$ws2.Cells.Item(15,12).Value = ""

# New row 16
$ws2.Cells.Item(16,1).Value = "Breakfast Oats"
$ws2.Cells.Item(16,2).Value = "yes"
$ws2.Cells.Item(16,3).Value = ""
$ws2.Cells.Item(16,4).Value = ""
$ws2.Cells.Item(16,5).Value = ""
$ws2.Cells.Item(16,6).Value = ""
$ws2.Cells.Item(16,7).Value = ""
$ws2.Cells.Item(16,8).Value = ""
$ws2.Cells.Item(16,9).Value = ""
$ws2.Cells.Item(16,10).Value = ""
$ws2.Cells.Item(16,11).Value = ""
$ws2.Cells.Item(16,12).Value = ""
